# Add OEE machine data to the configuration sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = "Machine 1"
$ws.Range("B4").Value = 600

$ws.Range("A5").Value = "Machine 2"
$ws.Range("B5").Value = 1200

$ws.Range("A6").Value = "Machine 3"
$ws.Range("B6").Value = 1500

$ws.Range("C6").Select()
